$wb = $excel.ActiveWorkbook

# --- Sheet 1 ("chức năng chính"): remove the three blank rows (16:18) ---
$ws1 = $wb.Worksheets.Item(1)
$ws1.Range("16:18").Delete() | Out-Null

# --- Sheet 2 ("thiết kế DB"): add the database schema documentation ---
$ws2 = $wb.Worksheets.Item(2)

# Table: SINH_VIEN (left column B, top to bottom)
$ws2.Range("B2").Value = "SINH_VIEN"
$ws2.Range("B3").Value = "ma"
$ws2.Range("B4").Value = "ten"
$ws2.Range("B5").Value = "gioi_tinh"
$ws2.Range("B6").Value = "cmnd"

# Table: TAI_KHOAN (column E, top to bottom)
$ws2.Range("E2").Value = "TAI_KHOAN"
$ws2.Range("E3").Value = "ma"
$ws2.Range("E4").Value = "mat_khau"
$ws2.Range("E5").Value = "liscen_id"

# Table: SINH_VIEN_LOP (header + first two fields)
$ws2.Range("B8").Value = "SINH_VIEN_LOP"
$ws2.Range("E8").Value = ""
$ws2.Range("B9").Value = "ma_sinh_vien"
$ws2.Range("B10").Value = "ma_lop"

# Table: LOP (column B)
$ws2.Range("B18").Value = "LOP"

# Table: MON_HOC (column H, bottom to top)
$ws2.Range("H21").Value = "phong_hoc"
$ws2.Range("H20").Value = "ten"
$ws2.Range("H19").Value = "ma"
$ws2.Range("H18").Value = "MON_HOC"

# Table: LOP_MON_HOC (column E)
$ws2.Range("E21").Value = "ngay_hoc"
$ws2.Range("E22").Value = "tiet_bat_dau"
$ws2.Range("E23").Value = "tiet_ket_thuc"
$ws2.Range("E18").Value = "LOP_MON_HOC"
$ws2.Range("E20").Value = "ma_mon_hoc"
$ws2.Range("E19").Value = "ma_lop"

# remaining LOP table fields
$ws2.Range("B19").Value = "ma"
$ws2.Range("B20").Value = "ten"

# Remaining SINH_VIEN_LOP fields (filled in after the LOP/MON_HOC tables)
$ws2.Range("B11").Value = "trang_thai_sinh_vien"
$ws2.Range("B12").Value = "diem_GK"
$ws2.Range("B13").Value = "diem_CK"
$ws2.Range("B14").Value = "diem_khac"
$ws2.Range("B15").Value = "diem_tong"

# Header shading (green) for table-name cells
$headerColor = 5296274  # RGB(146, 208, 80) == FF92D050
$ws2.Range("B2").Interior.Color = $headerColor
$ws2.Range("E2").Interior.Color = $headerColor
$ws2.Range("B8").Interior.Color = $headerColor
$ws2.Range("E8").Interior.Color = $headerColor
$ws2.Range("B18").Interior.Color = $headerColor
$ws2.Range("E18").Interior.Color = $headerColor
$ws2.Range("H18").Interior.Color = $headerColor

# Column widths to fit the longer header/labels
$ws2.Columns.Item(2).ColumnWidth = 19
$ws2.Columns.Item(5).ColumnWidth = 25.166666666666668

# --- Selection / active view state ---
$ws1.Range("B18").Select() | Out-Null

$ws2.Activate()
$ws2.Range("E11").Select() | Out-Null
